$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("accelerators template")
$wsv = $wb.Worksheets.Item("ValidationData")

function Set-ValidationFormula($sheet, $sqref, $formula) {
    $rng = $sheet.Range($sqref)
    $v = $rng.Validation
    $v.Modify($v.Type, $v.AlertStyle, $v.Operator, $formula)
}

$NValues = @("Software", "Fintech", "Healthcare", "Medtech", "AI", "Computing", "Deep tech", "Climate", "Consumer", "E-commerce", "Marketplace", "Gaming", "Web3", "Developer tools", "Cybersecurity", "Logistics", "Adtech", "Proptech", "Agriculture", "Automotive", "Biotechnology", "Construction", "Education", "Energy", "Entertainment", "Environment", "Fashion", "Real estate", "Food", "IoT", "Government", "Hospitality", "HR", "Insurance", "Security", "Social", "Aerospace", "AR/VR", "Mining", "Advanced Materials", "Biofuels", "Hardware", "Nanotechnology", "Legal", "Manufacturing", "Media", "Pharmaceuticals", "Retail", "Telecommunications", "Transportation", "Agnostic", "Other")
$OValues = @("Global", "North America", "South America", "LATAM", "Europe", "Middle East", "Africa", "Asia", "East Asia", "South East Asia", "South Asia", "Oceania", "EMEA", "Emerging Markets", "India", "China", "Japan", "Korea", "Australia", "United States", "Canada", "UK", "France", "Nigeria", "Kenya", "Egypt", "Senegal", "South Africa", "Netherlands", "Sweden", "Other")
$TValues = @("pitch_deck", "video", "financials", "business_plan")
$YValues = @("contact", "airtable", "typeform", "google", "generic")
$AAValues = @("FREE", "PRO", "MAX", "ENTERPRISE")

# Update column N (45 -> 52 items)
for ($i = 0; $i -lt $NValues.Length; $i++) {
    $wsv.Cells.Item($i + 1, 14).Value = $NValues[$i]
}

# Update column O (17 -> 31 items)
for ($i = 0; $i -lt $OValues.Length; $i++) {
    $wsv.Cells.Item($i + 1, 15).Value = $OValues[$i]
}

# Update column T (5 -> 4 items)
for ($i = 0; $i -lt $TValues.Length; $i++) {
    $wsv.Cells.Item($i + 1, 20).Value = $TValues[$i]
}
$wsv.Cells.Item(5, 20).ClearContents()

# Update column Y (4 -> 5 items)
for ($i = 0; $i -lt $YValues.Length; $i++) {
    $wsv.Cells.Item($i + 1, 25).Value = $YValues[$i]
}

# Update column AA (3 -> 4 items)
for ($i = 0; $i -lt $AAValues.Length; $i++) {
    $wsv.Cells.Item($i + 1, 27).Value = $AAValues[$i]
}

# Update data validation formulas for resized lists
Set-ValidationFormula $ws "N10:N1000" "ValidationData!`$N`$1:`$N`$52"
Set-ValidationFormula $ws "N2:N1000" "ValidationData!`$N`$1:`$N`$52"
Set-ValidationFormula $ws "O10:O1000" "ValidationData!`$O`$1:`$O`$31"
Set-ValidationFormula $ws "O2:O1000" "ValidationData!`$O`$1:`$O`$31"
Set-ValidationFormula $ws "T10:T1000" "ValidationData!`$T`$1:`$T`$4"
Set-ValidationFormula $ws "T2:T1000" "ValidationData!`$T`$1:`$T`$4"
Set-ValidationFormula $ws "Y10:Y1000" "ValidationData!`$Y`$1:`$Y`$5"
Set-ValidationFormula $ws "Y2:Y1000" "ValidationData!`$Y`$1:`$Y`$5"
Set-ValidationFormula $ws "AA10:AA1000" "ValidationData!`$AA`$1:`$AA`$4"
Set-ValidationFormula $ws "AA2:AA1000" "ValidationData!`$AA`$1:`$AA`$4"

# Update comments to reflect new possible values
$NComment = @"
Possible values:

- Software
- Fintech
- Healthcare
- Medtech
- AI
- Computing
- Deep tech
- Climate
- Consumer
- E-commerce
- Marketplace
- Gaming
- Web3
- Developer tools
- Cybersecurity
- Logistics
- Adtech
- Proptech
- Agriculture
- Automotive
- Biotechnology
- Construction
- Education
- Energy
- Entertainment
- Environment
- Fashion
- Real estate
- Food
- IoT
- Government
- Hospitality
- HR
- Insurance
- Security
- Social
- Aerospace
- AR/VR
- Mining
- Advanced Materials
- Biofuels
- Hardware
- Nanotechnology
- Legal
- Manufacturing
- Media
- Pharmaceuticals
- Retail
- Telecommunications
- Transportation
- Agnostic
- Other
"@
$ws.Range("N1").Comment.Text($NComment)

$OComment = @"
Possible values:

- Global
- North America
- South America
- LATAM
- Europe
- Middle East
- Africa
- Asia
- East Asia
- South East Asia
- South Asia
- Oceania
- EMEA
- Emerging Markets
- India
- China
- Japan
- Korea
- Australia
- United States
- Canada
- UK
- France
- Nigeria
- Kenya
- Egypt
- Senegal
- South Africa
- Netherlands
- Sweden
- Other
"@
$ws.Range("O1").Comment.Text($OComment)

$TComment = @"
Possible values:

- pitch_deck
- video
- financials
- business_plan
"@
$ws.Range("T1").Comment.Text($TComment)

$YComment = @"
Possible values:

- contact
- airtable
- typeform
- google
- generic
"@
$ws.Range("Y1").Comment.Text($YComment)

$AAComment = @"
Possible values:

- FREE
- PRO
- MAX
- ENTERPRISE
"@
$ws.Range("AA1").Comment.Text($AAComment)
